# Add Columbus Day (federal holiday, Mon Oct 11 2021) styling to the
# "1-15" sign-in sheet: narrow the V/W (Tech / Time of Arrival) columns
# the same way the existing Sat/Sun pairs are narrowed, mark them with
# an "X" (for the rows that already show "X" in the other off-day
# columns), and re-style the header/date/label rows to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1) Narrow columns V (22) and W (23) from 4.5 -> 2.5 characters.
#    (Other "weekend" pair columns, e.g. D/E, already use ColumnWidth 1.67
#    which round-trips to a stored width of 2.5.)
$ws.Range("V1").ColumnWidth = 1.67
$ws.Range("W1").ColumnWidth = 1.67

# Rows where the parallel Sat/Sun "off day" columns already carry an "X"
# (i.e. Columbus Day should likewise be marked as a non-working day).
$xRows = @(5, 6, 8, 9, 11, 12, 14, 15, 17, 18, 20, 21, 23, 24, 26, 27)

for ($row = 2; $row -le 27; $row++) {
    $vCell = $ws.Range("V$row")
    $wCell = $ws.Range("W$row")

    if ($xRows -contains $row) {
        $vCell.Value = "X"
        $wCell.Value = "X"
    }

    # Copy the cell formatting (borders/fill/font) from the D/E "SAT"
    # column pair, which already has the narrow off-day styling that
    # row $row needs for V/W - this leaves each cell's own value intact
    # and only updates its style.
    $ws.Range("D$row").Copy()
    $vCell.PasteSpecial(-4122)

    $ws.Range("E$row").Copy()
    $wCell.PasteSpecial(-4122)
}

$excel.CutCopyMode = 0
